# Sheet "Table1" (the active/selected tab) holds a list of YouTube short
# links in A2:A5, each cell carrying an external hyperlink. This upload
# clears those four link cells (keeping their style) and removes the
# hyperlinks entirely, then leaves the selection parked on R12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table1")

# Remove the link text from A2:A5 (keeps the existing cell style s="1").
$ws.Range("A2:A5").ClearContents()

# Drop the two external hyperlinks that used to live on A2 and A5.
$ws.Hyperlinks.Delete()

# Leave the selection where the author last clicked before saving.
$ws.Range("R12").Select() | Out-Null
